$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The multiple-comparison "contrast" column (e.g. "Winter / Summer") is split
# into two separate columns, "con_1" and "con_2" (e.g. "Winter", "Summer"),
# and a new "adj_p_value" column is appended after "statistic". All numeric
# results are refreshed per the updated model criteria.

# Header row
$ws.Cells.Item(1,1).Value = "term"
$ws.Cells.Item(1,2).Value = "con_1"
$ws.Cells.Item(1,3).Value = "con_2"
$ws.Cells.Item(1,4).Value = "null_value"
$ws.Cells.Item(1,5).Value = "ratio"
$ws.Cells.Item(1,6).Value = "std_error"
$ws.Cells.Item(1,7).Value = "df"
$ws.Cells.Item(1,8).Value = "null"
$ws.Cells.Item(1,9).Value = "statistic"
$ws.Cells.Item(1,10).Value = "adj_p_value"
# Row 2
$ws.Cells.Item(2,1).Value = "season"
$ws.Cells.Item(2,2).Value = "Winter"
$ws.Cells.Item(2,3).Value = "Summer"
$ws.Cells.Item(2,4).Value = 0
$ws.Cells.Item(2,5).Value = 0.884072295129919
$ws.Cells.Item(2,6).Value = 0.00356512345522638
$ws.Cells.Item(2,7).Value = "#NUM!"
$ws.Cells.Item(2,8).Value = 1
$ws.Cells.Item(2,9).Value = -30.554969664054
$ws.Cells.Item(2,10).Value = 0.00000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000291434732920951
# Row 3
$ws.Cells.Item(3,1).Value = "season"
$ws.Cells.Item(3,2).Value = "Winter"
$ws.Cells.Item(3,3).Value = "Spring"
$ws.Cells.Item(3,4).Value = 0
$ws.Cells.Item(3,5).Value = 0.883970328969684
$ws.Cells.Item(3,6).Value = 0.0044045850436028
$ws.Cells.Item(3,7).Value = "#NUM!"
$ws.Cells.Item(3,8).Value = 1
$ws.Cells.Item(3,9).Value = -24.7518516083878
$ws.Cells.Item(3,10).Value = 0.0000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000177671405641428
# Row 4
$ws.Cells.Item(4,1).Value = "season"
$ws.Cells.Item(4,2).Value = "Fall"
$ws.Cells.Item(4,3).Value = "Summer"
$ws.Cells.Item(4,4).Value = 0
$ws.Cells.Item(4,5).Value = 0.933263067255638
$ws.Cells.Item(4,6).Value = 0.00464176844058609
$ws.Cells.Item(4,7).Value = "#NUM!"
$ws.Cells.Item(4,8).Value = 1
$ws.Cells.Item(4,9).Value = -13.8866820093212
$ws.Cells.Item(4,10).Value = 0.000000000000000000000000000000000000000000457729480608557
# Row 5
$ws.Cells.Item(5,1).Value = "season"
$ws.Cells.Item(5,2).Value = "Fall"
$ws.Cells.Item(5,3).Value = "Winter"
$ws.Cells.Item(5,4).Value = 0
$ws.Cells.Item(5,5).Value = 1.0556411193934
$ws.Cells.Item(5,6).Value = 0.0046200193232886
$ws.Cells.Item(5,7).Value = "#NUM!"
$ws.Cells.Item(5,8).Value = 1
$ws.Cells.Item(5,9).Value = 12.3724913880289
$ws.Cells.Item(5,10).Value = 0.000000000000000000000000000000000220932036121536
# Row 6
$ws.Cells.Item(6,1).Value = "season"
$ws.Cells.Item(6,2).Value = "Fall"
$ws.Cells.Item(6,3).Value = "Spring"
$ws.Cells.Item(6,4).Value = 0
$ws.Cells.Item(6,5).Value = 0.933155427584107
$ws.Cells.Item(6,6).Value = 0.00544490255209406
$ws.Cells.Item(6,7).Value = "#NUM!"
$ws.Cells.Item(6,8).Value = 1
$ws.Cells.Item(6,9).Value = -11.8567707442956
$ws.Cells.Item(6,10).Value = 0.000000000000000000000000000000119090046741024
# Row 7
$ws.Cells.Item(7,1).Value = "season"
$ws.Cells.Item(7,2).Value = "Spring"
$ws.Cells.Item(7,3).Value = "Summer"
$ws.Cells.Item(7,4).Value = 0
$ws.Cells.Item(7,5).Value = 1.00011535020678
$ws.Cells.Item(7,6).Value = 0.00540875409542651
$ws.Cells.Item(7,7).Value = "#NUM!"
$ws.Cells.Item(7,8).Value = 1
$ws.Cells.Item(7,9).Value = 0.0213278062425114
$ws.Cells.Item(7,10).Value = 1
